$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.735.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.678.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +5.84%  "
$ws.Range("E9").Value = "  +5.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000199"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.160.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.615.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.678.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.31%  "
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("E28").Value = "  -5.88%  "
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "529.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "164.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0611"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.00%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.91%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0259"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.643"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("E51").Value = "  -4.30%  "
